$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Label the previously-empty header cell above the mutant/sample names
$ws.Range("A1").Value = "variant"

# Rename the first and last sample rows to reflect their biological role
# (first replicate set is the wild-type control, last is the triple mutant)
$ws.Range("A2").Value = "WT"
$ws.Range("A9").Value = "Triple"
